$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3808.6667
$ws.Range("I40").Value = 1682.8334
$ws.Range("J40").Value = 5934.5
$ws.Range("K40").Value = 1682.8334
$ws.Range("L40").Value = 5934.5
$ws.Range("M40").Value = -1507.8334
$ws.Range("N40").Value = -6284.5
$ws.Range("H64").Value = 5121.4
$ws.Range("I64").Value = 3381.2
$ws.Range("J64").Value = 6861.6
$ws.Range("K64").Value = 3381.2
$ws.Range("L64").Value = 6861.6
$ws.Range("M64").Value = -3133.2
$ws.Range("N64").Value = -7357.6
$ws.Range("H67").Value = 5121.4
$ws.Range("I67").Value = 3381.2
$ws.Range("J67").Value = 6861.6
$ws.Range("K67").Value = 3381.2
$ws.Range("L67").Value = 6861.6
$ws.Range("M67").Value = -2523.2
$ws.Range("N67").Value = -8577.6
$ws.Range("H93").Value = 133333
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 133333
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 133333
$ws.Range("N93").Value = -138325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 10789
$ws.Range("I46").Value = 28068.5
$ws.Range("J46").Value = 5029.1665
$ws.Range("K46").Value = 28068.5
$ws.Range("L46").Value = 5029.1665
$ws.Range("M46").Value = -27749.5
$ws.Range("N46").Value = -5667.1665
$ws.Range("H61").Value = 2646.2727
$ws.Range("I61").Value = 2630.9
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 2630.9
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -2418.9
$ws.Range("N61").Value = -3224
$ws.Range("H74").Value = 2885.8
$ws.Range("I74").Value = 3056.8823
$ws.Range("J74").Value = 1916.3334
$ws.Range("K74").Value = 3056.8823
$ws.Range("L74").Value = 1916.3334
$ws.Range("M74").Value = -2182.8823
$ws.Range("N74").Value = -3664.3334
$ws.Range("H77").Value = 2885.8
$ws.Range("I77").Value = 3056.8823
$ws.Range("J77").Value = 1916.3334
$ws.Range("K77").Value = 15284.4115
$ws.Range("L77").Value = 9581.666999999999
$ws.Range("M77").Value = -10916.4115
$ws.Range("N77").Value = -18317.667
$ws.Range("H110").Value = 18239.25
$ws.Range("I110").Value = 18239.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 18239.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -16194.25
$ws.Range("H131").Value = 79993
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 79993
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 79993
$ws.Range("N131").Value = -90073
$ws.Range("H132").Value = 3443.9546
$ws.Range("I132").Value = 3477.1191
$ws.Range("J132").Value = 2747.5
$ws.Range("K132").Value = 10431.3573
$ws.Range("L132").Value = 8242.5
$ws.Range("M132").Value = -7901.3573
$ws.Range("N132").Value = -13302.5
$ws.Range("H136").Value = 2646.2727
$ws.Range("I136").Value = 2630.9
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 7892.700000000001
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -5342.700000000001
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H94").Value = 1140.25
$ws.Range("I94").Value = 998.4
$ws.Range("J94").Value = 1849.5
$ws.Range("K94").Value = 998.4
$ws.Range("L94").Value = 1849.5
$ws.Range("M94").Value = -547.4
$ws.Range("N94").Value = -2751.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 402.25
$ws.Range("I2").Value = 203
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 203
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -90
$ws.Range("N2").Value = -1226
$ws.Range("H4").Value = 1483
$ws.Range("I4").Value = 1250
$ws.Range("H20").Value = 69994
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 69994
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 69994
$ws.Range("N20").Value = -70466
$ws.Range("H22").Value = 572.7857
$ws.Range("I22").Value = 634.7273
$ws.Range("J22").Value = 345.66666
$ws.Range("K22").Value = 634.7273
$ws.Range("L22").Value = 345.66666
$ws.Range("M22").Value = -284.7273
$ws.Range("N22").Value = -1045.66666
$ws.Range("N25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H30").Value = 69994
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 69994
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 69994
$ws.Range("N30").Value = -70176
$ws.Range("H58").Value = 3107.625
$ws.Range("I58").Value = 3194.4285
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 3194.4285
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -2991.4285
$ws.Range("N58").Value = -2906
$ws.Range("H128").Value = 69994
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 69994
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 69994
$ws.Range("N128").Value = -79954
$ws.Range("H132").Value = 2803.0833
$ws.Range("I132").Value = 2767
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 8301
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -5771
$ws.Range("N132").Value = -14660
$ws.Range("H134").Value = 5954298
$ws.Range("I134").Value = 6212963
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 18638889
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -18636354
$ws.Range("N134").Value = -20067
$ws.Range("H136").Value = 3107.625
$ws.Range("I136").Value = 3194.4285
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 9583.2855
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -7033.2855
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 927.25
$ws.Range("I19").Value = 772.3333
$ws.Range("J19").Value = 1392
$ws.Range("K19").Value = 2316.9999
$ws.Range("L19").Value = 4176
$ws.Range("M19").Value = -2142.9999
$ws.Range("N19").Value = -4524
$ws.Range("H70").Value = 13397.5
$ws.Range("I70").Value = 13397.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 40192.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -39877.5
$ws.Range("H73").Value = 13397.5
$ws.Range("I73").Value = 13397.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 40192.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -39100.5
$ws.Range("N99").ClearContents()
$ws.Range("H99").Value = 1499
$ws.Range("I99").Value = 1499
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4497
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2848.8667
$ws.Range("I122").Value = 2520.3333
$ws.Range("J122").Value = 3615.4443
$ws.Range("K122").Value = 7560.999899999999
$ws.Range("L122").Value = 10846.3329
$ws.Range("M122").Value = -5110.999899999999
$ws.Range("N122").Value = -15746.3329
$ws.Range("H132").Value = 3887.6667
$ws.Range("I132").Value = 3182.5557
$ws.Range("J132").Value = 6003
$ws.Range("K132").Value = 9547.667099999999
$ws.Range("L132").Value = 18009
$ws.Range("M132").Value = -7017.667099999999
$ws.Range("N132").Value = -23069

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 999
$ws.Range("I16").Value = 998
$ws.Range("J16").Value = 1001
$ws.Range("K16").Value = 998
$ws.Range("L16").Value = 1001
$ws.Range("M16").Value = -828
$ws.Range("N16").Value = -1341
$ws.Range("H132").Value = 3624.2222
$ws.Range("I132").Value = 4483.2856
$ws.Range("J132").Value = 3077.5454
$ws.Range("K132").Value = 13449.8568
$ws.Range("L132").Value = 9232.636200000001
$ws.Range("M132").Value = -10919.8568
$ws.Range("N132").Value = -14292.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 13499.5
$ws.Range("I56").Value = 15000
$ws.Range("J56").Value = 11999
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 11999
$ws.Range("M56").Value = -14286
$ws.Range("N56").Value = -13427
$ws.Range("H80").Value = 25998.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 25998.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25998.5
$ws.Range("N80").Value = -27994.5
$ws.Range("H83").Value = 25998.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 25998.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 77995.5
$ws.Range("N83").Value = -87979.5
$ws.Range("H92").Value = 13997.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 13997.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 13997.5
$ws.Range("N92").Value = -18989.5
$ws.Range("H132").Value = 3616.7666
$ws.Range("I132").Value = 2473.3044
$ws.Range("J132").Value = 7373.857
$ws.Range("K132").Value = 7419.9132
$ws.Range("L132").Value = 22121.571
$ws.Range("M132").Value = -4889.9132
$ws.Range("N132").Value = -27181.571
